$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.909.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.455.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.510'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.455.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('E13').Value = '  -5.51%  '
$ws.Range('E15').Value = '  -5.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.924.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('E17').Value = '  -5.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.453.29'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.90%  '
$ws.Range('E20').Value = '  -8.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.63%  '
$ws.Range('E22').Value = '  -4.18%  '
$ws.Range('E23').Value = '  -1.98%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.86%  '
$ws.Range('E26').Value = '  -8.44%  '
$ws.Range('E27').Value = '  -5.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.63%  '
$ws.Range('E29').Value = '  -30.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.582.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '507.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0886'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.85%  '
$ws.Range('E33').Value = '  -8.95%  '
$ws.Range('E34').Value = '  -6.33%  '
$ws.Range('E35').Value = '  -7.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('E38').Value = '  -12.55%  '
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.30%  '
$ws.Range('E41').Value = '  -8.87%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.45%  '
$ws.Range('E44').Value = '  -7.42%  '
$ws.Range('E45').Value = '  -7.95%  '
$ws.Range('E46').Value = '  -6.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '140.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.56%  '
$ws.Range('E49').Value = '  -8.08%  '
$ws.Range('E50').Value = '  -9.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0726'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.16%  '
